$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "Terms Typically Offered" column (D),
# shifting it to column G. This makes room for the new Corequisites (D),
# Concurrent (E) and Recommended (F) columns.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row labels
$ws.Cells.Item(1, 4).Value2 = "Corequisites"
$ws.Cells.Item(1, 5).Value2 = "Concurrent"
$ws.Cells.Item(1, 6).Value2 = "Recommended"

# Default every data row (2-20) in the new D, E, F columns to "NA"
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "NA"
    $ws.Cells.Item($r, 5).Value2 = "NA"
    $ws.Cells.Item($r, 6).Value2 = "NA"
}

# Row 6 (AGC 301) previously embedded a "Recommended:" clause inside the
# Prerequisites text. Split it out: trim the Prerequisites cell and move the
# recommendation into the new Recommended column.
$ws.Cells.Item(6, 3).Value2 = "Junior standing."
$ws.Cells.Item(6, 6).Value2 = "JOUR 203, JOUR 205."

# The "Terms Typically Offered" value for row 6 gains a trailing space.
$ws.Cells.Item(6, 7).Value2 = "W "

# Make sure the sheet dimension reflects the new data bounds.
$ws.Range("A1:G20").Select()
